$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.691.98"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.001.58"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "379.81"
$ws.Range("E5").Value = "  +4.46%  "
$ws.Range("D6").Value = "105.45"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "37.84"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "18.82"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "3.471.71"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "7.55"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "2.996.59"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "0.965"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "51.630.41"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "3.51"
$ws.Range("E19").Value = "  +5.57%  "
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").Value = "13.20"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").Value = "68.91"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "264.54"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "2.79"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").Value = "7.36"
$ws.Range("E26").Value = "  +17.16%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "26.17"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").Value = "9.97"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "34.58"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").Value = "51.46"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "0.0444"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").Value = "17.55"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("E40").Value = "  -5.97%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Value = "124.23"
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "22.59"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").Value = "0.283"
$ws.Range("E45").Value = "  +18.81%  "
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").Value = "2.46"
$ws.Range("E47").Value = "  +7.33%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.050.81"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "3.30"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "0.0349"
$ws.Range("E50").Value = "  +9.48%  "
$ws.Range("D51").Value = "0.875"
$ws.Range("E51").Value = "  -0.39%  "
